# Populate the cross-validation comparison table (rows 2-8, columns A-P)
# with the new multi-condition best-per-run summary data.
#
# Cells are written column-by-column (A down to P, top-to-bottom within each
# column) so that newly introduced shared strings are appended to
# xl/sharedStrings.xml in the same order as the canonical export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "run_g3"
$ws.Range("A3").Value = "run_g4"
$ws.Range("A4").Value = "run_g5"
$ws.Range("A5").Value = "run_g2"
$ws.Range("A6").Value = "run_g1"
$ws.Range("A7").Value = "run_g6"
$ws.Range("A8").Value = "TOTAL"
$ws.Range("B2").Value = "Multi-Condition"
$ws.Range("B3").Value = "Multi-Condition"
$ws.Range("B4").Value = "Multi-Condition"
$ws.Range("B5").Value = "Multi-Condition"
$ws.Range("B6").Value = "Multi-Condition"
$ws.Range("B7").Value = "Multi-Condition"
$ws.Range("B8").Value = "TOTAL"
$ws.Range("C2").Value = "double_AND"
$ws.Range("C3").Value = "double_AND"
$ws.Range("C4").Value = "double_AND"
$ws.Range("C5").Value = "double_AND"
$ws.Range("C6").Value = "double_AND"
$ws.Range("C7").Value = "double_AND"
$ws.Range("C8").Value = "TOTAL"
$ws.Range("D2").Value = "anom_swh_min_swan > t1 AND swh_max_swan > t2"
$ws.Range("D3").Value = "swh_p80_swan > t1 AND anom_swh_mean_deseasonalized_detrended > t2"
$ws.Range("D4").Value = "swh_p80_swan > t1 AND swh_max_swan > t2"
$ws.Range("D5").Value = "anom_swh_min_waverys > t1 AND anom_swh_max_waverys > t2"
$ws.Range("D6").Value = "anom_swe_min > t1 AND anom_swh_max_waverys > t2"
$ws.Range("D7").Value = "anom_swh_min_waverys > t1 AND anom_swh_p80_swan > t2"
$ws.Range("D8").Value = "TOTAL"
$ws.Range("E2").Value = 0.3307086614173228
$ws.Range("E3").Value = 0.5521064301552105
$ws.Range("E4").Value = 0.4281524926686216
$ws.Range("E5").Value = 0.296875
$ws.Range("E6").Value = 0.4583333333333333
$ws.Range("E7").Value = 0.2681451612903225
$ws.Range("E8").Value = 2.334321078864811
$ws.Range("F2").Value = 1711241.429632653
$ws.Range("F3").Value = 5601398
$ws.Range("F4").Value = 570803.0660659341
$ws.Range("F5").Value = 243425.9018571429
$ws.Range("F6").Value = 81190.8
$ws.Range("F7").Value = 3237740.665735714
$ws.Range("F8").Value = 11445799.86329144
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0.015
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 0.015
$ws.Range("H2").Value = 21
$ws.Range("H3").Value = 249
$ws.Range("H4").Value = 73
$ws.Range("H5").Value = 19
$ws.Range("H6").Value = 22
$ws.Range("H7").Value = 133
$ws.Range("H8").Value = 517
$ws.Range("I2").Value = 46
$ws.Range("I3").Value = 269
$ws.Range("I4").Value = 55
$ws.Range("I5").Value = 56
$ws.Range("I6").Value = 27
$ws.Range("I7").Value = 589
$ws.Range("I8").Value = 1042
$ws.Range("J2").Value = 2245
$ws.Range("J3").Value = 1904
$ws.Range("J4").Value = 2289
$ws.Range("J5").Value = 2448
$ws.Range("J6").Value = 2483
$ws.Range("J7").Value = 1705
$ws.Range("J8").Value = 13074
$ws.Range("K2").Value = 29
$ws.Range("K3").Value = 135
$ws.Range("K4").Value = 140
$ws.Range("K5").Value = 34
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 137
$ws.Range("K8").Value = 500
$ws.Range("L2").Value = 0.3134328358208955
$ws.Range("L3").Value = 0.4806949806949807
$ws.Range("L4").Value = 0.5703125
$ws.Range("L5").Value = 0.2533333333333334
$ws.Range("L6").Value = 0.4489795918367347
$ws.Range("L7").Value = 0.1842105263157895
$ws.Range("L8").Value = 2.250963768001734
$ws.Range("M2").Value = 0.42
$ws.Range("M3").Value = 0.6484375
$ws.Range("M4").Value = 0.3427230046948357
$ws.Range("M5").Value = 0.3584905660377358
$ws.Range("M6").Value = 0.4680851063829787
$ws.Range("M7").Value = 0.4925925925925926
$ws.Range("M8").Value = 2.730328769708143
$ws.Range("N2").Value = 0.9679624092268262
$ws.Range("N3").Value = 0.8420023464998044
$ws.Range("N4").Value = 0.9237387563551036
$ws.Range("N5").Value = 0.9648025029331248
$ws.Range("N6").Value = 0.979663668361361
$ws.Range("N7").Value = 0.7168486739469578
$ws.Range("N8").Value = 5.395018357323178
$ws.Range("O2").Value = 0.358974358974359
$ws.Range("O3").Value = 0.5521064301552107
$ws.Range("O4").Value = 0.4281524926686217
$ws.Range("O5").Value = 0.296875
$ws.Range("O6").Value = 0.4583333333333333
$ws.Range("O7").Value = 0.2681451612903226
$ws.Range("O8").Value = 2.362586776421847
$ws.Range("P2").Value = "enhanced_multi_rule_summary_20250711_143735.csv"
$ws.Range("P3").Value = "enhanced_multi_rule_summary_20250714_121146.csv"
$ws.Range("P4").Value = "enhanced_multi_rule_summary_20250714_124208.csv"
$ws.Range("P5").Value = "enhanced_multi_rule_summary_20250713_172447.csv"
$ws.Range("P6").Value = "enhanced_multi_rule_summary_20250713_170835.csv"
$ws.Range("P7").Value = "enhanced_multi_rule_summary_20250714_130045.csv"
$ws.Range("P8").Value = "TOTAL"

